$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "35.264.39"
$ws.Range("E2").Value = "  +1.23%  "
# Row 3
$ws.Range("D3").Value = "1.883.24"
$ws.Range("E3").Value = "  +0.36%  "
# Row 4
$ws.Range("E4").Value = "  -0.01%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.98"
$ws.Range("E5").Value = "  -0.64%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.685"
$ws.Range("E6").Value = "  +0.09%  "
# Row 7
$ws.Range("E7").Value = "  +0.07%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.60"
$ws.Range("E8").Value = "  +1.60%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.352"
$ws.Range("E9").Value = "  +1.90%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "55.08"
$ws.Range("E10").Value = "  +7.91%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0740"
$ws.Range("E11").Value = "  +0.49%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0982"
$ws.Range("E12").Value = "  +1.29%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "13.79"
$ws.Range("E13").Value = "  +7.33%  "
# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.159.21"
$ws.Range("E14").Value = "  +0.48%  "
# Row 15
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.771"
$ws.Range("E15").Value = "  +8.01%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.96"
$ws.Range("E16").Value = "  +2.10%  "
# Row 17
$ws.Range("D17").Value = "1.873.58"
$ws.Range("E17").Value = "  -0.05%  "
# Row 18
$ws.Range("D18").Value = "35.251.78"
$ws.Range("E18").Value = "  +1.16%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.12"
$ws.Range("E19").Value = "  +0.65%  "
# Row 20
$ws.Range("D20").Value = "0.0₃0819"
$ws.Range("E20").Value = "  +0.28%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "243.70"
$ws.Range("E21").Value = "  +0.12%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.74"
$ws.Range("E22").Value = "  +0.72%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.12"
$ws.Range("E23").Value = "  +4.33%  "
# Row 24
$ws.Range("E24").Value = "  +7.46%  "
# Row 26
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.05"
$ws.Range("E26").Value = "  +1.36%  "
# Row 27
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.13"
$ws.Range("E27").Value = "  -3.12%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.49"
$ws.Range("E28").Value = "  +1.49%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.20"
$ws.Range("E29").Value = "  +0.30%  "
# Row 30
$ws.Range("E30").Value = "  +0.28%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.27"
$ws.Range("E31").Value = "  +1.45%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0590"
$ws.Range("E32").Value = "  +2.27%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.89"
$ws.Range("E33").Value = "  +23.70%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.16"
$ws.Range("E34").Value = "  +0.50%  "
# Row 35
$ws.Range("E35").Value = "  -0.01%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.45"
$ws.Range("E36").Value = "  -13.66%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.844"
$ws.Range("E37").Value = "  +2.49%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.92"
$ws.Range("E38").Value = "  -2.69%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0719"
$ws.Range("E39").Value = "  +8.18%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0220"
$ws.Range("E40").Value = "  +4.13%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.41"
$ws.Range("E41").Value = "  +0.05%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.06"
$ws.Range("E42").Value = "  +0.57%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.06"
$ws.Range("E43").Value = "  -1.87%  "
# Row 44
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.327.22"
$ws.Range("E44").Value = "  +3.47%  "
# Row 45
$ws.Range("B45").Value = "Gas"
$ws.Range("C45").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.06"
$ws.Range("E45").Value = "  +10.81%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.36"
$ws.Range("E46").Value = "  +2.07%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0807"
$ws.Range("E47").Value = "  -0.46%  "
# Row 48
$ws.Range("E48").Value = "  +0.13%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.73"
$ws.Range("E49").Value = "  +0.20%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.26"
$ws.Range("E50").Value = "  -2.29%  "
# Row 51
$ws.Range("D51").Value = "2.058.85"
$ws.Range("E51").Value = "  +0.07%  "